$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the D1 header from "Cond (1 truth, 2 lie)" to "Condition"
$ws.Range("D1").Value = "Condition"

# Add a new "word count" column in E
$ws.Range("E1").Value = "word count"

$wordCounts = @(7,9,6,8,8,14,14,17,9,15,7,27,10,25,17,15,21,24,46,25,36,34,28,28,38,32,32,24,32,36,34,17,32,28,38,37,51,25,42,28,18,20,23,19,19,12,29,31,29,27,32,12,21,16,24,27,18,15,55,41,59,38,51,52,47,51,42,42,37,37,39,42,42,31,27,34,43,21,41,41,29,38,18,62,36,55,28,27,40,18,27,29,40,29,35,45,38,25,79,49,74,64,71,11,53,66,43,64,30,45,37,38,51,35,48,46,35,30,64,29,23,43,23,26,23,24,18,18,28,26,28,38,28,29,33,29,54,50,50,26,28,36,31,47,36,34,61,23,42,43,33,30,27,24,30,30,28,31,22,26,29,18,29,20,28,30,19,22,29,22,22,22,28,21,28,29,29,31,22,23,31,39,27,24,18,15,32,33,32,25,30,22,37,18,26,24,26,28,37,37,31,32,38,34,41,39,25,34,36,28,24,34,21,35,20,31,24,27,30,30,37,42,30,26,26,14,31,17,10,11,19,11,15,20,27,22,24,30,28,20,23,30,31,23,21,28,31,15,22,47,39,45,84,52,48,62,17,27,16,6,34,27)
for ($i = 0; $i -lt $wordCounts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $wordCounts[$i]
}

# Widen column D slightly (matches the author's manual resize)
$ws.Columns.Item(4).ColumnWidth = 10

# Move the active selection to G10 (matches the saved view state)
$ws.Range("G10").Select() | Out-Null
